$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.940.72"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "2.340.46"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "306.79"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").Value = "100.92"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  -4.97%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -3.81%  "
$ws.Range("D10").Value = "34.91"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("D11").Value = "52.04"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D14").Value = "6.80"
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "15.79"
$ws.Range("E15").Value = "  +5.22%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.806"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.244.62"
$ws.Range("E17").Value = "  -3.24%  "
$ws.Range("D18").Value = "42.859.20"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").Value = "11.69"
$ws.Range("E21").Value = "  -5.96%  "
$ws.Range("D22").Value = "67.85"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "236.93"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("D24").Value = "2.01"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "25.41"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("E28").Value = "  +9.41%  "
$ws.Range("D29").Value = "35.05"
$ws.Range("E29").Value = "  -4.62%  "
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("D31").Value = "160.12"
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("D34").Value = "4.68"
$ws.Range("E34").Value = "  +9.74%  "
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "17.43"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.0728"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.024.67"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0286"
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "18.68"
$ws.Range("E44").Value = "  -4.15%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "10.31"
$ws.Range("E45").Value = "  +3.65%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "2.95"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("B47").Value = "MultiversX"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D47").Value = "56.34"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "2.90"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.566.16"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "4.66"
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.52"
$ws.Range("E51").Value = "  -3.82%  "
